# Applies cryptos.xlsx price/volume updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.508.39"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "1.831.52"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4291"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3652"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07276"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8683"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "1.933.92"
$ws.Range("E12").Value = "  +7.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.408"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.532"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06931"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "80.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008889"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "27.486.81"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.132"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.18%  "
$ws.Range("D24").Value = "2.090.32"
$ws.Range("E24").Value = "  +3.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.979"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  -1.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.145"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.832"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08885"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7557"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.988"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.537"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.133"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05311"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01937"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.799"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1661"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5071"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.604"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.376"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "105.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06498"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4680"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.609"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.19%  "
